# Corrections to the training-data matrix on sheet "Matriz"
# (mirrors fixes made in neurona.py's logic).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: only B2 changes (2 -> 1)
$ws.Range("B2").Value = 1

# Row 3: A3 1.5 -> 1, B3 2 -> 0, C3 1 -> 0
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0

# Row 4: only A4 changes (1 -> 0)
$ws.Range("A4").Value = 0

# Row 5: A5 2 -> 0, B5 2 -> 0 (and B5 loses its special border/alignment style)
$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B5").ClearFormats()

# Update the sheet's active selection to F14
$ws.Range("F14").Select()
